# Oregon overview factsheet - text edits from COMM
#
# The "No. of 990 Filers w/ Gov Grants" counts were stored as real numbers.
# They need to become plain text values instead (so e.g. 1961 -> "1,961"
# with a thousands separator baked into the text itself). A new "Total"
# row is also added at the bottom of the County sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper data: for every sheet, the column-B counts that must become text
# (keyed by row number), plus the label used for the "Total" row (if any).
# ---------------------------------------------------------------------

$wsOverall = $wb.Worksheets.Item("Overall")
$wsCounty  = $wb.Worksheets.Item("County")
$wsCD      = $wb.Worksheets.Item("Congressional District")
$wsSize    = $wb.Worksheets.Item("Size")
$wsSub     = $wb.Worksheets.Item("Subsector")

# --- Overall sheet: A2 is the lone filer count ------------------------
$wsOverall.Range("A2").NumberFormat = "@"
$wsOverall.Range("A2").Value = "1,961"

# --- County sheet: B2:B37 counts -> text, plus new Total row 38 -------
$countyCounts = @{
  2="11"; 3="63"; 4="114"; 5="39"; 6="18"; 7="25"; 8="14"; 9="11"; 10="99";
  11="33"; 12="2"; 13="6"; 14="6"; 15="15"; 16="107"; 17="3"; 18="30";
  19="34"; 20="9"; 21="183"; 22="25"; 23="37"; 24="15"; 25="108"; 26="6";
  27="636"; 28="22"; 29="1"; 30="19"; 31="28"; 32="14"; 33="9"; 34="18";
  35="169"; 36="1"; 37="31"
}

$wsCounty.Range("B2:B37").NumberFormat = "@"
foreach ($r in $countyCounts.Keys) {
  $wsCounty.Cells.Item($r, 2).Value = $countyCounts[$r]
}

$wsCounty.Range("A38").Value = "Total"
$wsCounty.Range("B38:F38").NumberFormat = "@"
$wsCounty.Range("B38").Value = "1,961"
$wsCounty.Range("C38").Value = "$3,508,813,229"
$wsCounty.Range("D38").Value = "11.60%"
$wsCounty.Range("E38").Value = "-13.87%"
$wsCounty.Range("F38").Value = "67.31%"

# --- Congressional District sheet: B2:B7 counts -> text, B8 Total -----
$cdCounts = @{ 2="510"; 3="318"; 4="328"; 5="340"; 6="256"; 7="209" }

$wsCD.Range("B2:B8").NumberFormat = "@"
foreach ($r in $cdCounts.Keys) {
  $wsCD.Cells.Item($r, 2).Value = $cdCounts[$r]
}
$wsCD.Range("B8").Value = "1,961"

# --- Size sheet: B2:B7 counts -> text, B8 Total ------------------------
$sizeCounts = @{ 2="631"; 3="561"; 4="368"; 5="140"; 6="188"; 7="73" }

$wsSize.Range("B2:B8").NumberFormat = "@"
foreach ($r in $sizeCounts.Keys) {
  $wsSize.Cells.Item($r, 2).Value = $sizeCounts[$r]
}
$wsSize.Range("B8").Value = "1,961"

# --- Subsector sheet: B2:B12 counts -> text, B13 Total -----------------
$subCounts = @{
  2="177"; 3="232"; 4="172"; 5="153"; 6="11"; 7="622"; 8="17"; 9="147";
  10="41"; 11="376"; 12="13"
}

$wsSub.Range("B2:B13").NumberFormat = "@"
foreach ($r in $subCounts.Keys) {
  $wsSub.Cells.Item($r, 2).Value = $subCounts[$r]
}
$wsSub.Range("B13").Value = "1,961"
